$wb = $excel.ActiveWorkbook
$dataSheet = $wb.Worksheets.Item(1)

# --- Update the panel_query_time timestamps in the "data" sheet (column F) ---
$newTimes = @("2021-10-05 14:19:27.757887","2021-10-05 14:19:27.757898","2021-10-05 14:19:27.757902","2021-10-05 14:19:27.757905","2021-10-05 14:19:27.757908","2021-10-05 14:19:27.757910","2021-10-05 14:19:27.757913","2021-10-05 14:19:27.757916","2021-10-05 14:19:27.757919","2021-10-05 14:19:27.757921","2021-10-05 14:19:27.757924","2021-10-05 14:19:27.757927","2021-10-05 14:19:27.757929","2021-10-05 14:19:27.757932","2021-10-05 14:19:27.757935","2021-10-05 14:19:27.757937","2021-10-05 14:19:27.757940","2021-10-05 14:19:27.757943","2021-10-05 14:19:27.757946","2021-10-05 14:19:27.757949","2021-10-05 14:19:27.757952","2021-10-05 14:19:27.757954","2021-10-05 14:19:27.757957","2021-10-05 14:19:27.757960","2021-10-05 14:19:27.757963","2021-10-05 14:19:27.757966","2021-10-05 14:19:27.757969","2021-10-05 14:19:27.757972","2021-10-05 14:19:27.757974","2021-10-05 14:19:27.757977","2021-10-05 14:19:27.757980","2021-10-05 14:19:27.757982","2021-10-05 14:19:27.757985","2021-10-05 14:19:27.757989","2021-10-05 14:19:27.757991","2021-10-05 14:19:27.757994","2021-10-05 14:19:27.757997","2021-10-05 14:19:27.757999","2021-10-05 14:19:27.758002","2021-10-05 14:19:27.758005","2021-10-05 14:19:27.758008","2021-10-05 14:19:27.758010","2021-10-05 14:19:27.758013","2021-10-05 14:19:27.758016","2021-10-05 14:19:27.758019","2021-10-05 14:19:27.758021","2021-10-05 14:19:27.758024","2021-10-05 14:19:27.758027","2021-10-05 14:19:27.758030","2021-10-05 14:19:27.758032","2021-10-05 14:19:27.758035","2021-10-05 14:19:27.758038","2021-10-05 14:19:27.758041","2021-10-05 14:19:27.758044","2021-10-05 14:19:27.758046","2021-10-05 14:19:27.758049","2021-10-05 14:19:27.758052","2021-10-05 14:19:27.758054","2021-10-05 14:19:27.758057","2021-10-05 14:19:27.758060","2021-10-05 14:19:27.758062","2021-10-05 14:19:27.758065","2021-10-05 14:19:27.758067","2021-10-05 14:19:27.758070","2021-10-05 14:19:27.758074","2021-10-05 14:19:27.758077","2021-10-05 14:19:27.758080","2021-10-05 14:19:27.758083","2021-10-05 14:19:27.758085","2021-10-05 14:19:27.758088","2021-10-05 14:19:27.758090","2021-10-05 14:19:27.758093","2021-10-05 14:19:27.758096","2021-10-05 14:19:27.758099","2021-10-05 14:19:27.758101")
for ($i = 0; $i -lt $newTimes.Length; $i++) {
    $dataSheet.Cells.Item($i + 2, 6).Value = $newTimes[$i]
}

# --- Add a new "metadata" worksheet right after "data" ---
$metaSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $dataSheet)
$metaSheet.Name = "metadata"

# Header row (row 1), columns B..G, bold/border/centered via the same style
# used for the "data" sheet header (style index 1 == Range.Style "header"-like)
$metaSheet.Cells.Item(1, 2).Value = "data_name"
$metaSheet.Cells.Item(1, 3).Value = "data_id"
$metaSheet.Cells.Item(1, 4).Value = "data_version"
$metaSheet.Cells.Item(1, 5).Value = "data_version_created"
$metaSheet.Cells.Item(1, 6).Value = "panel_query_time"
$metaSheet.Cells.Item(1, 7).Value = "panel_get_request"

$headerRange = $metaSheet.Range("B1:G1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160
$headerRange.Borders.LineStyle = 1

# Data row (row 2)
$metaSheet.Cells.Item(2, 1).Value = 0
$metaSheet.Cells.Item(2, 1).Font.Bold = $true
$metaSheet.Cells.Item(2, 1).HorizontalAlignment = -4108
$metaSheet.Cells.Item(2, 1).VerticalAlignment = -4160
$metaSheet.Cells.Item(2, 1).Borders.LineStyle = 1

$metaSheet.Cells.Item(2, 2).Value = "Cerebellar hypoplasia"
$metaSheet.Cells.Item(2, 3).Value = 286

$metaSheet.Range("D2").NumberFormat = "@"
$metaSheet.Cells.Item(2, 4).Value = "1.59"

$metaSheet.Cells.Item(2, 5).Value = "2021-09-01T13:40:25.074429Z"
$metaSheet.Cells.Item(2, 6).Value = "2021-10-05 14:19:27.754346"
$metaSheet.Cells.Item(2, 7).Value = "https://panelapp.genomicsengland.co.uk/api/v1/panels/286/?format=json"

$dataSheet.Select()
Write-Output "ok"
